$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1040.3636
$ws.Range("I43").Value = 795
$ws.Range("J43").Value = 1180.5714
$ws.Range("K43").Value = 795
$ws.Range("L43").Value = 1180.5714
$ws.Range("M43").Value = -726
$ws.Range("N43").Value = -1318.5714
$ws.Range("H116").Value = 4570.294
$ws.Range("I116").Value = 2699.5454
$ws.Range("J116").Value = 8000
$ws.Range("K116").Value = 2699.5454
$ws.Range("L116").Value = 8000
$ws.Range("M116").Value = 742.4546
$ws.Range("N116").Value = -14884
$ws.Range("H132").Value = 4819.7915
$ws.Range("I132").Value = 4040.9424
$ws.Range("K132").Value = 12122.8272
$ws.Range("M132").Value = -9592.8272
$ws.Range("H138").Value = 1671.0588
$ws.Range("I138").Value = 796
$ws.Range("J138").Value = 7171.4287
$ws.Range("K138").Value = 2388
$ws.Range("L138").Value = 21514.2861
$ws.Range("M138").Value = 2752
$ws.Range("N138").Value = -31794.2861

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2421.32
$ws.Range("I45").Value = 1196.9231
$ws.Range("J45").Value = 3747.75
$ws.Range("K45").Value = 1196.9231
$ws.Range("L45").Value = 3747.75
$ws.Range("M45").Value = -819.9231
$ws.Range("N45").Value = -4501.75
$ws.Range("H102").Value = 2000
$ws.Range("I102").Value = 2000
$ws.Range("K102").Value = 2000
$ws.Range("M102").Value = -378
$ws.Range("H132").Value = 3818.7036
$ws.Range("I132").Value = 1587.2333
$ws.Range("J132").Value = 6608.0415
$ws.Range("K132").Value = 4761.699900000001
$ws.Range("L132").Value = 19824.1245
$ws.Range("M132").Value = -2231.699900000001
$ws.Range("N132").Value = -24884.1245

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1495.2858
$ws.Range("I107").Value = 1789.5834
$ws.Range("J107").Value = 1102.8889
$ws.Range("K107").Value = 1789.5834
$ws.Range("L107").Value = 1102.8889
$ws.Range("M107").Value = 130.4166
$ws.Range("N107").Value = -4942.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2759.4443
$ws.Range("I16").Value = 2830.9167
$ws.Range("J16").Value = 2616.5
$ws.Range("K16").Value = 2830.9167
$ws.Range("L16").Value = 2616.5
$ws.Range("M16").Value = -2543.9167
$ws.Range("N16").Value = -3190.5
$ws.Range("H99").Value = 2771.68
$ws.Range("I99").Value = 2633.7144
$ws.Range("K99").Value = 2633.7144
$ws.Range("M99").Value = -1135.7144
$ws.Range("H107").Value = 908.7083
$ws.Range("I107").Value = 348.5
$ws.Range("J107").Value = 1693
$ws.Range("K107").Value = 348.5
$ws.Range("L107").Value = 1693
$ws.Range("M107").Value = 1571.5
$ws.Range("N107").Value = -5533
$ws.Range("H113").Value = 2759.4443
$ws.Range("I113").Value = 2830.9167
$ws.Range("J113").Value = 2616.5
$ws.Range("K113").Value = 2830.9167
$ws.Range("L113").Value = 2616.5
$ws.Range("M113").Value = -660.9167000000002
$ws.Range("N113").Value = -6956.5
$ws.Range("H122").Value = 71430250
$ws.Range("I122").Value = 125000840
$ws.Range("J122").Value = 2799.8333
$ws.Range("K122").Value = 375002520
$ws.Range("L122").Value = 8399.499899999999
$ws.Range("M122").Value = -375000070
$ws.Range("N122").Value = -13299.4999
$ws.Range("H126").Value = 2771.68
$ws.Range("I126").Value = 2633.7144
$ws.Range("K126").Value = 2633.7144
$ws.Range("L126").Value = 7901.1432
$ws.Range("M126").Value = -5431.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 242.88461
$ws.Range("I2").Value = 41.63158
$ws.Range("J2").Value = 789.1429000000001
$ws.Range("K2").Value = 249.78948
$ws.Range("L2").Value = 4734.857400000001
$ws.Range("M2").Value = -136.78948
$ws.Range("N2").Value = -4960.857400000001
$ws.Range("H92").Value = 1204.9231
$ws.Range("I92").Value = 770.25
$ws.Range("J92").Value = 1900.4
$ws.Range("K92").Value = 2310.75
$ws.Range("L92").Value = 5701.200000000001
$ws.Range("M92").Value = -1062.75
$ws.Range("N92").Value = -8197.200000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3972058.2
$ws.Range("I102").Value = 11908528
$ws.Range("J102").Value = 3823.5
$ws.Range("K102").Value = 11908528
$ws.Range("L102").Value = 3823.5
$ws.Range("M102").Value = -11906906
$ws.Range("N102").Value = -7067.5
$ws.Range("H122").Value = 372284.38
$ws.Range("I122").Value = 436607.4
$ws.Range("J122").Value = 2427
$ws.Range("K122").Value = 1309822.2
$ws.Range("L122").Value = 7281
$ws.Range("M122").Value = -1307372.2
$ws.Range("N122").Value = -12181
$ws.Range("H132").Value = 852431.4
$ws.Range("I132").Value = 1603997.1
$ws.Range("J132").Value = 2835.3044
$ws.Range("K132").Value = 4811991.300000001
$ws.Range("L132").Value = 8505.913199999999
$ws.Range("M132").Value = -4809461.300000001
$ws.Range("N132").Value = -13565.9132

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 780.93335
$ws.Range("I22").Value = 225
$ws.Range("J22").Value = 983.0909
$ws.Range("K22").Value = 225
$ws.Range("L22").Value = 983.0909
$ws.Range("M22").Value = 70
$ws.Range("N22").Value = -1573.0909
$ws.Range("H27").Value = 780.93335
$ws.Range("I27").Value = 225
$ws.Range("J27").Value = 983.0909
$ws.Range("K27").Value = 225
$ws.Range("L27").Value = 983.0909
$ws.Range("M27").Value = -118
$ws.Range("N27").Value = -1197.0909
$ws.Range("H55").Value = 345.11765
$ws.Range("I55").Value = 268.91666
$ws.Range("J55").Value = 528
$ws.Range("K55").Value = 268.91666
$ws.Range("L55").Value = 528
$ws.Range("M55").Value = -95.91665999999998
$ws.Range("N55").Value = -874
$ws.Range("H61").Value = 4751.273
$ws.Range("I61").Value = 1612.8
$ws.Range("J61").Value = 7366.6665
$ws.Range("K61").Value = 1612.8
$ws.Range("L61").Value = 7366.6665
$ws.Range("M61").Value = -1410.8
$ws.Range("N61").Value = -7770.6665
$ws.Range("H113").Value = 4751.273
$ws.Range("I113").Value = 1612.8
$ws.Range("J113").Value = 7366.6665
$ws.Range("K113").Value = 1612.8
$ws.Range("L113").Value = 7366.6665
$ws.Range("M113").Value = 557.2
$ws.Range("N113").Value = -11706.6665
$ws.Range("H136").Value = 17859404
$ws.Range("I136").Value = 27779012
$ws.Range("J136").Value = 4110.5
$ws.Range("K136").Value = 83337036
$ws.Range("L136").Value = 12331.5
$ws.Range("M136").Value = -83334486
$ws.Range("N136").Value = -17431.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3710.147
$ws.Range("I122").Value = 2923.7036
$ws.Range("J122").Value = 6743.5713
$ws.Range("K122").Value = 8771.110799999999
$ws.Range("L122").Value = 20230.7139
$ws.Range("M122").Value = -6321.110799999999
$ws.Range("N122").Value = -25130.7139
$ws.Range("H136").Value = 5210088
$ws.Range("I136").Value = 7353971.5
$ws.Range("J136").Value = 3513.9285
$ws.Range("K136").Value = 22061914.5
$ws.Range("L136").Value = 10541.7855
$ws.Range("M136").Value = -22059364.5
$ws.Range("N136").Value = -15641.7855
